$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty formatted row 22 (collapses to an absent row, matching target)
$ws.Rows.Item(22).Clear()

# Remove the old task "6" row (Pendiente / placeholder data) - replaced by new tasks below
$ws.Rows.Item(33).Clear()

# Template formatting (style s=2 "Salida" cells, s=3 date cells) copied from row 2
$ws.Range("A2:G2").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)
$ws.Range("A30:G30").PasteSpecial(-4122)
$ws.Range("A31:G31").PasteSpecial(-4122)
$ws.Range("A32:G32").PasteSpecial(-4122)
$ws.Range("A33:G33").PasteSpecial(-4122)
$ws.Range("A34:G34").PasteSpecial(-4122)
$ws.Range("A35:G35").PasteSpecial(-4122)
$ws.Range("A36:G36").PasteSpecial(-4122)
$ws.Range("A37:G37").PasteSpecial(-4122)
$ws.Range("A38:G38").PasteSpecial(-4122)
$ws.Range("A39:G39").PasteSpecial(-4122)
$ws.Range("A40:G40").PasteSpecial(-4122)
$ws.Range("A41:G41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A29").Value2 = 18
$ws.Range("B29").Value2 = 'Añadir botón de eliminación y edición dentro del detalle'
$ws.Range("C29").Value2 = 'Jon'
$ws.Range("D29").Value2 = 45795
$ws.Range("E29").Value2 = 45795
$ws.Range("F29").Value2 = '✅ Hecho'
$ws.Range("G29").Value2 = ""
$ws.Rows.Item(29).RowHeight = 30

$ws.Range("A30").Value2 = 19
$ws.Range("B30").Value2 = 'Correcciones anterior entrega - Añadir CRUD restantes a Clientes Productos y Componentes'
$ws.Range("C30").Value2 = 'Jon'
$ws.Range("D30").Value2 = 45795
$ws.Range("E30").Value2 = 45795
$ws.Range("F30").Value2 = '✅ Hecho'
$ws.Range("G30").Value2 = ""
$ws.Rows.Item(30).RowHeight = 45

$ws.Range("A31").Value2 = 20
$ws.Range("B31").Value2 = 'Correcciones anterior entrega - Añadir widgets de tipo fecha en los formularios de creacion'
$ws.Range("C31").Value2 = 'Jon'
$ws.Range("D31").Value2 = 45796
$ws.Range("E31").Value2 = 45796
$ws.Range("F31").Value2 = '✅ Hecho'
$ws.Range("G31").Value2 = ""
$ws.Rows.Item(31).RowHeight = 45

$ws.Range("A32").Value2 = 21
$ws.Range("B32").Value2 = 'Ampliación de funcionalidades en Python - Envió de emails'
$ws.Range("C32").Value2 = 'Jon'
$ws.Range("D32").Value2 = 45796
$ws.Range("E32").Value2 = 45798
$ws.Range("F32").Value2 = '✅ Hecho'
$ws.Range("G32").Value2 = ""
$ws.Rows.Item(32).RowHeight = 30

$ws.Range("A33").Value2 = 22
$ws.Range("B33").Value2 = 'Ampliación de funcionalidades en Python - Subida de ficheros al servidor mediante <input type="file"> y mostrarlos en una página (tienen que poder descargarse)'
$ws.Range("C33").Value2 = 'Xabier'
$ws.Range("D33").Value2 = 45797
$ws.Range("E33").Value2 = 45798
$ws.Range("F33").Value2 = '✅ Hecho'
$ws.Range("G33").Value2 = ""
$ws.Rows.Item(33).RowHeight = 75

$ws.Range("A34").Value2 = 23
$ws.Range("B34").Value2 = 'Implementaciones JS - Implementar las siguientes funcionalidades JavaScript Capturar un evento en el DOM y producir un cambio en el estilo/ mostrar una alerta si el usuario realiza una acción determinada,...)'
$ws.Range("C34").Value2 = 'Jon'
$ws.Range("D34").Value2 = 45799
$ws.Range("E34").Value2 = 45800
$ws.Range("F34").Value2 = '✅ Hecho'
$ws.Range("G34").Value2 = ""
$ws.Rows.Item(34).RowHeight = 90

$ws.Range("A35").Value2 = 24
$ws.Range("B35").Value2 = 'Cambios Django - Mejorar el sistema de precios del sistema, en vez de meter el precio directamente en el producto, hacerlo en componentes, para luego así poder implementar el #7 Calcular autom.'
$ws.Range("C35").Value2 = 'Xabier'
$ws.Range("D35").Value2 = 45799
$ws.Range("E35").Value2 = 45800
$ws.Range("F35").Value2 = '✅ Hecho'
$ws.Range("G35").Value2 = ""
$ws.Rows.Item(35).RowHeight = 75

$ws.Range("A36").Value2 = 25
$ws.Range("B36").Value2 = 'Implementaciones JS - Validar campos de un formulario antes de su envío al servidor'
$ws.Range("C36").Value2 = 'Jon'
$ws.Range("D36").Value2 = 45799
$ws.Range("E36").Value2 = 45800
$ws.Range("F36").Value2 = '✅ Hecho'
$ws.Range("G36").Value2 = ""
$ws.Rows.Item(36).RowHeight = 45

$ws.Range("A37").Value2 = 26
$ws.Range("B37").Value2 = 'Implementaciones JS - Calcular automáticamente el precio total de un pedido utilizando API y Fetch'
$ws.Range("C37").Value2 = 'Xabier'
$ws.Range("D37").Value2 = 45800
$ws.Range("E37").Value2 = 45801
$ws.Range("F37").Value2 = '✅ Hecho'
$ws.Range("G37").Value2 = ""
$ws.Rows.Item(37).RowHeight = 45

$ws.Range("A38").Value2 = 26
$ws.Range("B38").Value2 = 'Implementaciones JS - Convertir caracteres seleccionados a mayusculas'
$ws.Range("C38").Value2 = 'Xabier'
$ws.Range("D38").Value2 = 45801
$ws.Range("E38").Value2 = 45803
$ws.Range("F38").Value2 = '✅ Hecho'
$ws.Range("G38").Value2 = ""
$ws.Rows.Item(38).RowHeight = 30

$ws.Range("A39").Value2 = 27
$ws.Range("B39").Value2 = 'Ampliacion de funcionalidades en Python - Paginación en tablas/listados de los resultados de una tabla.'
$ws.Range("C39").Value2 = 'Jon'
$ws.Range("D39").Value2 = 45802
$ws.Range("E39").Value2 = 45802
$ws.Range("F39").Value2 = '✅ Hecho'
$ws.Range("G39").Value2 = ""
$ws.Rows.Item(39).RowHeight = 45

$ws.Range("A40").Value2 = 28
$ws.Range("B40").Value2 = 'Implementaciones JS - Rango deslizante para controlar el tamaño del texto en tiempo real.'
$ws.Range("C40").Value2 = 'Jon'
$ws.Range("D40").Value2 = 45802
$ws.Range("E40").Value2 = 45803
$ws.Range("F40").Value2 = '✅ Hecho'
$ws.Range("G40").Value2 = ""
$ws.Rows.Item(40).RowHeight = 45

$ws.Range("A41").Value2 = 29
$ws.Range("B41").Value2 = 'Realizacion de los videos y ultimos retoques'
$ws.Range("C41").Value2 = 'Jon, Xabier'
$ws.Range("D41").Value2 = 45803
$ws.Range("E41").Value2 = 45804
$ws.Range("F41").Value2 = '✅ Hecho'
$ws.Range("G41").Value2 = ""
$ws.Rows.Item(41).RowHeight = 30

# Row 15 height adjustment (60 -> 45)
$ws.Rows.Item(15).RowHeight = 45

# View state (scroll position / active selection) to match the saved workbook state
$ws.Application.ActiveWindow.ScrollRow = 53
$ws.Range("E43").Select()

Write-Output "done"